$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 1111111111
$ws.Range("B1").Value = 9811126101
$ws.Range("A2").Value = 2222222222
$ws.Range("B2").Value = 9811126034
$ws.Range("A3").Value = 3333333333
$ws.Range("B3").Value = 9811126041
$ws.Range("A4").Value = 1234567890
$ws.Range("B4").Value = 9811126033
$ws.Range("A5").Value = 2345678910
$ws.Range("B5").Value = 9811126072
$ws.Range("A6").Value = 3456789120
$ws.Range("B6").Value = 9811126087

$ws.Range("A7:B20").Value = 0
$ws.Range("A7:B20").ClearContents()

$ws.Range("B2").Select()
